# Update cryptocurrency price/volume data pulled from coinranking.com
# (refresh performed by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '54.254.08'
$ws.Cells.Item(2, 5).Value = '  +0.05%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.263.03'
$ws.Cells.Item(3, 5).Value = '  -1.03%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.30%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''495.76'
$ws.Cells.Item(5, 5).Value = '  -0.05%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''128.63'
$ws.Cells.Item(6, 5).Value = '  +0.38%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.01%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.79%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.50%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.91%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.336'
$ws.Cells.Item(11, 5).Value = '  +2.60%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''4.80'
$ws.Cells.Item(12, 5).Value = '  +3.39%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''22.94'
$ws.Cells.Item(13, 5).Value = '  +5.06%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.663.38'
$ws.Cells.Item(14, 5).Value = '  -0.98%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '54.234.79'
$ws.Cells.Item(15, 5).Value = '  +0.05%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +0.12%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.274.81'
$ws.Cells.Item(17, 5).Value = '  -0.54%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.86%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.58%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''302.55'
$ws.Cells.Item(20, 5).Value = '  +0.34%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -2.08%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''0.998'
$ws.Cells.Item(22, 5).Value = '  -0.12%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''60.71'
$ws.Cells.Item(23, 5).Value = '  -2.83%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -1.36%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.01%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +3.37%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''171.23'
$ws.Cells.Item(27, 5).Value = '  +1.27%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''1.60'
$ws.Cells.Item(28, 5).Value = '  -0.06%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''5.96'
$ws.Cells.Item(29, 5).Value = '  +1.76%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.03%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +1.01%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.00%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''17.77'
$ws.Cells.Item(33, 5).Value = '  +0.44%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.12%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.936'
$ws.Cells.Item(35, 5).Value = '  +6.99%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.11%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.21%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.01%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''1.40'
$ws.Cells.Item(39, 5).Value = '  -0.31%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''3.36'
$ws.Cells.Item(40, 5).Value = '  +0.06%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(41, 4).Value = '''4.81'
$ws.Cells.Item(41, 5).Value = '  +0.48%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(42, 4).Value = '''124.68'
$ws.Cells.Item(42, 5).Value = '  -2.16%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.0493'
$ws.Cells.Item(43, 5).Value = '  +1.95%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.64%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''0.545'
$ws.Cells.Item(45, 5).Value = '  +0.21%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''241.45'
$ws.Cells.Item(46, 5).Value = '  +1.28%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.05%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.82%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.99%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''16.11'
$ws.Cells.Item(50, 5).Value = '  -1.05%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.41%  '
